# Update "想去人数" (interest count) values in the "展览" and "全部类型" sheets.
# 展览 (sheet 1): F2 578->582, F3 124->125, F4 27->29, F6 347->348, F7 1532->1577
# 全部类型 (sheet 4): same edits, but the last row lives at F11 instead of F7.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 582
$wsExhibition.Range("F3").Value = 125
$wsExhibition.Range("F4").Value = 29
$wsExhibition.Range("F6").Value = 348
$wsExhibition.Range("F7").Value = 1577

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 582
$wsAll.Range("F3").Value = 125
$wsAll.Range("F4").Value = 29
$wsAll.Range("F6").Value = 348
$wsAll.Range("F11").Value = 1577
